$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------------
# New "fired" column (I) - values first, in the order that makes the new
# shared-string table read: toto, --, fired (24, 25, 26)
# ------------------------------------------------------------------------

# New row 8 content (values) -- B8/C8 introduce the new shared strings in
# the right order, then I1 introduces "fired" last.
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "toto"
$ws.Range("C8").Value = "'--"
$ws.Range("D8").Value = 1
$ws.Range("F8").Value = "KO"
$ws.Range("G8").Value = "'false"
$ws.Range("H8").Value = "N/A"

$ws.Range("I1").Value = "fired"
$ws.Range("I2").Value = 1
$ws.Range("I3").Value = 1
$ws.Range("I4").Value = 1
$ws.Range("I5").Value = 1
$ws.Range("I6").Value = 1
$ws.Range("I7").Value = 1
$ws.Range("I8").Value = 0

# ------------------------------------------------------------------------
# Formats: copy the per-row look from the matching existing cells so the
# new cells line up with the table's existing banding. (G8 is re-pasted
# after its value so the quote-prefix picked up from the leading "'" used
# to force text instead of a boolean is cleared again.)
# ------------------------------------------------------------------------

$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)

$ws.Range("G2").Copy()
$ws.Range("I2").PasteSpecial(-4122)

$ws.Range("G3").Copy()
$ws.Range("I3").PasteSpecial(-4122)

$ws.Range("G4").Copy()
$ws.Range("I4").PasteSpecial(-4122)

$ws.Range("G5").Copy()
$ws.Range("I5").PasteSpecial(-4122)

$ws.Range("G6").Copy()
$ws.Range("I6").PasteSpecial(-4122)

$ws.Range("G7").Copy()
$ws.Range("I7").PasteSpecial(-4122)

$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)

$ws.Range("B7").Copy()
$ws.Range("B8").PasteSpecial(-4122)

$ws.Range("D7").Copy()
$ws.Range("D8").PasteSpecial(-4122)

$ws.Range("E7").Copy()
$ws.Range("E8").PasteSpecial(-4122)

$ws.Range("F7").Copy()
$ws.Range("F8").PasteSpecial(-4122)

$ws.Range("G7").Copy()
$ws.Range("G8").PasteSpecial(-4122)

$ws.Range("H7").Copy()
$ws.Range("H8").PasteSpecial(-4122)

$ws.Range("G7").Copy()
$ws.Range("I8").PasteSpecial(-4122)

# ------------------------------------------------------------------------
# View: the author landed on the new text cell after typing the row.
# ------------------------------------------------------------------------
$ws.Range("C8").Select()
